$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append two new rows (22 and 23) of daily data for 2026-01-11 (Excel serial 46033),
# one per charging station, continuing the existing table.
# Clone formatting/styles from the prior two rows (20:21, a matching
# station pair) so the new rows reuse the existing cellXfs (no bloat in styles.xml).
$ws.Range("A20:F21").Copy()
$ws.Range("A22:F23").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 22: 四方坪站 station
$ws.Range("A22").Value = 46033
$ws.Range("B22").Value = "四方坪站"
$ws.Range("C22").Value = 13186.83
$ws.Range("D22").Value = 9019.52
$ws.Range("E22").Value = 2918.73
$ws.Range("F22").Value = 541

# Row 23: 高岭站 station
$ws.Range("A23").Value = 46033
$ws.Range("B23").Value = "高岭站"
$ws.Range("C23").Value = 4627.6899999999996
$ws.Range("D23").Value = 3892.22
$ws.Range("E23").Value = 1313.1
$ws.Range("F23").Value = 176

# Update the selection to match the diff
$ws.Range("I20").Select()
